{"js": "// Requirement-list edit: add \"\uc870\ud68c \uc694\uccad \uc2dc\" lead-ins to requirements #6 and #7,\n// add a \"delete button\" mention to #6, and add the parenthetical detail list to #7.\nconst body = context.document.body;\n\nasync function replaceOnce(searchText, replaceText) {\n  const results = body.search(searchText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      `Expected exactly one match for \"${searchText}\", found ${results.items.length}`\n    );\n  }\n\n  results.items[0].insertText(replaceText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// ---- Requirement #6 (\ub300\uc5ec\uc18c \ub9ac\uc2a4\ud2b8 \uc870\ud68c) ----\n// \"\uc2dc\uc2a4\ud15c\uc740 \uad00\ub9ac\uc790\uac00 \ub4f1\ub85d\ud55c \ub300\uc5ec\uc18c \ub9ac\uc2a4\ud2b8\ub97c \uc81c\uacf5\ud55c\ub2e4.\"\n//   -> \"\uc2dc\uc2a4\ud15c\uc740 \uad00\ub9ac\uc790\uc758 \ub300\uc5ec\uc18c \ub9ac\uc2a4\ud2b8 \uc870\ud68c \uc694\uccad \uc2dc, \uad00\ub9ac\uc790\uac00 \ub4f1\ub85d\ud55c \ub300\uc5ec\uc18c \ub9ac\uc2a4\ud2b8\uc640 \uac01 \ud56d\ubaa9\uc5d0 \uc0ad\uc81c \ubc84\ud2bc\uc744 \uc81c\uacf5\ud55c\ub2e4. \"\n// Use a unique, longer anchor (\"\uad00\ub9ac\uc790\uac00 \ub4f1\ub85d\ud55c\") so this does not collide with\n// the other occurrences of the bare word \"\uad00\ub9ac\uc790\" elsewhere in the table.\nawait replaceOnce(\"\uad00\ub9ac\uc790\uac00 \ub4f1\ub85d\ud55c\", \"\uad00\ub9ac\uc790\uc758 \ub300\uc5ec\uc18c \ub9ac\uc2a4\ud2b8 \uc870\ud68c \uc694\uccad \uc2dc, \uad00\ub9ac\uc790\uac00 \ub4f1\ub85d\ud55c\");\nawait replaceOnce(\"\ub300\uc5ec\uc18c \ub9ac\uc2a4\ud2b8\ub97c \uc81c\uacf5\ud55c\ub2e4.\", \"\ub300\uc5ec\uc18c \ub9ac\uc2a4\ud2b8\uc640 \uac01 \ud56d\ubaa9\uc5d0 \uc0ad\uc81c \ubc84\ud2bc\uc744 \uc81c\uacf5\ud55c\ub2e4. \");\n\n// ---- Requirement #7 (\ub300\uc5ec\uc18c \uc0c1\uc138 \uc815\ubcf4 \uc870\ud68c) ----\n// \"\uc2dc\uc2a4\ud15c\uc740 \uad00\ub9ac\uc790\uac00 \ud574\ub2f9 \ub300\uc5ec\uc18c\ub97c \ub4f1\ub85d\ud560 \ub54c \uc785\ub825\ud55c \uc0c1\uc138 \uc815\ubcf4\ub97c \uc81c\uacf5\ud55c\ub2e4.\"\n//   -> \"\uc2dc\uc2a4\ud15c\uc740 \uad00\ub9ac\uc790\uc758 \ub300\uc5ec\uc18c \uc0c1\uc138 \uc815\ubcf4 \uc870\ud68c \uc694\uccad \uc2dc, \uad00\ub9ac\uc790\uac00 \ud574\ub2f9 \ub300\uc5ec\uc18c\ub97c \ub4f1\ub85d\ud560 \ub54c\n//       \uc785\ub825\ud55c \uc0c1\uc138 \uc815\ubcf4(\ub300\uc5ec\uc18c \uc774\ub984, \uc704\uce58, \uc790\uc804\uac70 \ubcf4\uad00 \uac00\ub2a5 \uc218\ub7c9, \uc6b4\uc601 \uc2dc\uac04)\ub97c \uc81c\uacf5\ud55c\ub2e4.\"\nawait replaceOnce(\"\uad00\ub9ac\uc790\uac00 \ud574\ub2f9 \ub300\uc5ec\uc18c\", \"\uad00\ub9ac\uc790\uc758 \ub300\uc5ec\uc18c \uc0c1\uc138 \uc815\ubcf4 \uc870\ud68c \uc694\uccad \uc2dc, \uad00\ub9ac\uc790\uac00 \ud574\ub2f9 \ub300\uc5ec\uc18c\");\nawait replaceOnce(\n  \" \uc785\ub825\ud55c \uc0c1\uc138 \uc815\ubcf4\ub97c \uc81c\uacf5\ud55c\ub2e4.\",\n  \" \uc785\ub825\ud55c \uc0c1\uc138 \uc815\ubcf4(\ub300\uc5ec\uc18c \uc774\ub984, \uc704\uce58, \uc790\uc804\uac70 \ubcf4\uad00 \uac00\ub2a5 \uc218\ub7c9, \uc6b4\uc601 \uc2dc\uac04)\ub97c \uc81c\uacf5\ud55c\ub2e4.\"\n);\n", "ps1": "# Requirement-list edit: add \"\uc870\ud68c \uc694\uccad \uc2dc\" lead-ins to requirements #6 and #7,\n# add a \"delete button\" mention to #6, and add the parenthetical detail list to #7.\n$d = $word.ActiveDocument\n\nfunction Replace-OnceInDoc($searchText, $replaceText) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $searchText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $result = $find.Execute()\n    if (-not $result) {\n        throw \"Could not find text: $searchText\"\n    }\n    $range.Text = $replaceText\n}\n\n# ---- Requirement #6 (\ub300\uc5ec\uc18c \ub9ac\uc2a4\ud2b8 \uc870\ud68c) ----\n# \"\uc2dc\uc2a4\ud15c\uc740 \uad00\ub9ac\uc790\uac00 \ub4f1\ub85d\ud55c \ub300\uc5ec\uc18c \ub9ac\uc2a4\ud2b8\ub97c \uc81c\uacf5\ud55c\ub2e4.\"\n#   -> \"\uc2dc\uc2a4\ud15c\uc740 \uad00\ub9ac\uc790\uc758 \ub300\uc5ec\uc18c \ub9ac\uc2a4\ud2b8 \uc870\ud68c \uc694\uccad \uc2dc, \uad00\ub9ac\uc790\uac00 \ub4f1\ub85d\ud55c \ub300\uc5ec\uc18c \ub9ac\uc2a4\ud2b8\uc640 \uac01 \ud56d\ubaa9\uc5d0 \uc0ad\uc81c \ubc84\ud2bc\uc744 \uc81c\uacf5\ud55c\ub2e4. \"\nReplace-OnceInDoc \"\uad00\ub9ac\uc790\uac00 \ub4f1\ub85d\ud55c\" \"\uad00\ub9ac\uc790\uc758 \ub300\uc5ec\uc18c \ub9ac\uc2a4\ud2b8 \uc870\ud68c \uc694\uccad \uc2dc, \uad00\ub9ac\uc790\uac00 \ub4f1\ub85d\ud55c\"\nReplace-OnceInDoc \"\ub300\uc5ec\uc18c \ub9ac\uc2a4\ud2b8\ub97c \uc81c\uacf5\ud55c\ub2e4.\" \"\ub300\uc5ec\uc18c \ub9ac\uc2a4\ud2b8\uc640 \uac01 \ud56d\ubaa9\uc5d0 \uc0ad\uc81c \ubc84\ud2bc\uc744 \uc81c\uacf5\ud55c\ub2e4. \"\n\n# ---- Requirement #7 (\ub300\uc5ec\uc18c \uc0c1\uc138 \uc815\ubcf4 \uc870\ud68c) ----\n# \"\uc2dc\uc2a4\ud15c\uc740 \uad00\ub9ac\uc790\uac00 \ud574\ub2f9 \ub300\uc5ec\uc18c\ub97c \ub4f1\ub85d\ud560 \ub54c \uc785\ub825\ud55c \uc0c1\uc138 \uc815\ubcf4\ub97c \uc81c\uacf5\ud55c\ub2e4.\"\n#   -> \"\uc2dc\uc2a4\ud15c\uc740 \uad00\ub9ac\uc790\uc758 \ub300\uc5ec\uc18c \uc0c1\uc138 \uc815\ubcf4 \uc870\ud68c \uc694\uccad \uc2dc, \uad00\ub9ac\uc790\uac00 \ud574\ub2f9 \ub300\uc5ec\uc18c\ub97c \ub4f1\ub85d\ud560 \ub54c\n#       \uc785\ub825\ud55c \uc0c1\uc138 \uc815\ubcf4(\ub300\uc5ec\uc18c \uc774\ub984, \uc704\uce58, \uc790\uc804\uac70 \ubcf4\uad00 \uac00\ub2a5 \uc218\ub7c9, \uc6b4\uc601 \uc2dc\uac04)\ub97c \uc81c\uacf5\ud55c\ub2e4.\"\nReplace-OnceInDoc \"\uad00\ub9ac\uc790\uac00 \ud574\ub2f9 \ub300\uc5ec\uc18c\" \"\uad00\ub9ac\uc790\uc758 \ub300\uc5ec\uc18c \uc0c1\uc138 \uc815\ubcf4 \uc870\ud68c \uc694\uccad \uc2dc, \uad00\ub9ac\uc790\uac00 \ud574\ub2f9 \ub300\uc5ec\uc18c\"\nReplace-OnceInDoc \" \uc785\ub825\ud55c \uc0c1\uc138 \uc815\ubcf4\ub97c \uc81c\uacf5\ud55c\ub2e4.\" \" \uc785\ub825\ud55c \uc0c1\uc138 \uc815\ubcf4(\ub300\uc5ec\uc18c \uc774\ub984, \uc704\uce58, \uc790\uc804\uac70 \ubcf4\uad00 \uac00\ub2a5 \uc218\ub7c9, \uc6b4\uc601 \uc2dc\uac04)\ub97c \uc81c\uacf5\ud55c\ub2e4.\"\n"}
